$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2.869551889135607
$ws.Cells.Item(2, 3).Value = 0.7881453474972204
$ws.Cells.Item(2, 4).Value = 0.07747209130812394
$ws.Cells.Item(2, 5).Value = 0.01810567887071102
$ws.Cells.Item(2, 7).Value = 0.002517548526561747
$ws.Cells.Item(2, 13).Value = 0.7062854332487376
$ws.Cells.Item(2, 14).Value = 3.017181111994375

$ws.Cells.Item(3, 2).Value = 2.643780808096665
$ws.Cells.Item(3, 3).Value = 0.7127389179150896
$ws.Cells.Item(3, 4).Value = 0.07037155588928101
$ws.Cells.Item(3, 5).Value = 0.01670266240005702
$ws.Cells.Item(3, 7).Value = 0.002526523239203815
$ws.Cells.Item(3, 13).Value = 0.6478379289009126
$ws.Cells.Item(3, 14).Value = 2.926055007387163

$ws.Cells.Item(4, 2).Value = 2.507240872029115
$ws.Cells.Item(4, 3).Value = 0.6669794408791745
$ws.Cells.Item(4, 4).Value = 0.06606562063427646
$ws.Cells.Item(4, 5).Value = 0.01585501730446737
$ws.Cells.Item(4, 7).Value = 0.002532304644592566
$ws.Cells.Item(4, 13).Value = 0.6124627174711605
$ws.Cells.Item(4, 14).Value = 2.870474595041543

$ws.Cells.Item(5, 2).Value = 2.452109964015222
$ws.Cells.Item(5, 3).Value = 0.6484630465233749
$ws.Cells.Item(5, 4).Value = 0.06432400157021334
$ws.Cells.Item(5, 5).Value = 0.0155129683418096
$ws.Cells.Item(5, 7).Value = 0.002534729055474287
$ws.Cells.Item(5, 13).Value = 0.5981720084138402
$ws.Cells.Item(5, 14).Value = 2.847913872052231

$ws.Cells.Item(6, 2).Value = 2.442985933878219
$ws.Cells.Item(6, 3).Value = 0.6453961814277704
$ws.Cells.Item(6, 4).Value = 0.06403558435978596
$ws.Cells.Item(6, 5).Value = 0.01545637222508134
$ws.Cells.Item(6, 7).Value = 0.002535135770275143
$ws.Cells.Item(6, 13).Value = 0.5958064881790506
$ws.Cells.Item(6, 14).Value = 2.844172906754267

$ws.Cells.Item(7, 2).Value = 2.506495309620107
$ws.Cells.Item(7, 3).Value = 0.6667291984668395
$ws.Cells.Item(7, 4).Value = 0.06604208015932045
$ws.Cells.Item(7, 5).Value = 0.0158503907849088
$ws.Cells.Item(7, 7).Value = 0.002532337063485515
$ws.Cells.Item(7, 13).Value = 0.6122694868005709
$ws.Cells.Item(7, 14).Value = 2.870169980002913

$ws.Cells.Item(8, 2).Value = 2.791264862275852
$ws.Cells.Item(8, 3).Value = 0.7620300433865737
$ws.Cells.Item(8, 4).Value = 0.07501238048971004
$ws.Cells.Item(8, 5).Value = 0.01761899075170348
$ws.Cells.Item(8, 7).Value = 0.002520586978306554
$ws.Cells.Item(8, 13).Value = 0.6860242326414436
$ws.Cells.Item(8, 14).Value = 2.985680913108126

$ws.Cells.Item(9, 2).Value = 3.366905249674232
$ws.Cells.Item(9, 3).Value = 0.9534375110640099
$ws.Cells.Item(9, 4).Value = 0.09305169667854329
$ws.Cells.Item(9, 5).Value = 0.02120157692914404
$ws.Cells.Item(9, 7).Value = 0.002499679712230204
$ws.Cells.Item(9, 13).Value = 0.8348962041147274
$ws.Cells.Item(9, 14).Value = 3.215366421890849

$ws.Cells.Item(10, 2).Value = 3.801306469134943
$ws.Cells.Item(10, 3).Value = 1.097167145946855
$ws.Cells.Item(10, 4).Value = 0.1066103141328654
$ws.Cells.Item(10, 5).Value = 0.02391041044955244
$ws.Cells.Item(10, 7).Value = 0.002485599561745577
$ws.Cells.Item(10, 13).Value = 0.947119284627945
$ws.Cells.Item(10, 14).Value = 3.386374647276256

$ws.Cells.Item(11, 2).Value = 4.001642273623133
$ws.Cells.Item(11, 3).Value = 1.163305008333737
$ws.Cells.Item(11, 4).Value = 0.112851873727962
$ws.Cells.Item(11, 5).Value = 0.02516095789603412
$ws.Cells.Item(11, 7).Value = 0.002479467634781252
$ws.Cells.Item(11, 13).Value = 0.9988498725304709
$ws.Cells.Item(11, 14).Value = 3.464733961479567

$ws.Cells.Item(12, 2).Value = 4.07791393930745
$ws.Cells.Item(12, 3).Value = 1.188464475656758
$ws.Cells.Item(12, 4).Value = 0.115226558998188
$ws.Cells.Item(12, 5).Value = 0.02563726421093548
$ws.Cells.Item(12, 7).Value = 0.00247718456380948
$ws.Cells.Item(12, 13).Value = 1.018541415603778
$ws.Cells.Item(12, 14).Value = 3.49449384462622

$ws.Cells.Item(13, 2).Value = 4.061469010885162
$ws.Cells.Item(13, 3).Value = 1.183040744761058
$ws.Cells.Item(13, 4).Value = 0.114714623601742
$ws.Cells.Item(13, 5).Value = 0.02553455879157696
$ws.Cells.Item(13, 7).Value = 0.002477674536618965
$ws.Cells.Item(13, 13).Value = 1.014295866210361
$ws.Cells.Item(13, 14).Value = 3.488080568590874

$ws.Cells.Item(14, 2).Value = 4.007908896077538
$ws.Cells.Item(14, 3).Value = 1.165372563261144
$ws.Cells.Item(14, 4).Value = 0.1130470140164732
$ws.Cells.Item(14, 5).Value = 0.02520008805672447
$ws.Cells.Item(14, 7).Value = 0.002479279026311379
$ws.Cells.Item(14, 13).Value = 1.000467831031045
$ws.Cells.Item(14, 14).Value = 3.467180551776153

$ws.Cells.Item(15, 2).Value = 3.97515553784126
$ws.Cells.Item(15, 3).Value = 1.15456538330227
$ws.Cells.Item(15, 4).Value = 0.1120270225440976
$ws.Cells.Item(15, 5).Value = 0.02499557702413568
$ws.Cells.Item(15, 7).Value = 0.002480266886231911
$ws.Cells.Item(15, 13).Value = 0.9920112188313794
$ws.Cells.Item(15, 14).Value = 3.454390168825853

$ws.Cells.Item(16, 2).Value = 3.788270324525115
$ws.Cells.Item(16, 3).Value = 1.092860572635232
$ws.Cells.Item(16, 4).Value = 0.1062039425433738
$ws.Cells.Item(16, 5).Value = 0.02382906251026995
$ws.Cells.Item(16, 7).Value = 0.002486005766998334
$ws.Cells.Item(16, 13).Value = 0.9437526328458716
$ws.Cells.Item(16, 14).Value = 3.381265546880428

$ws.Cells.Item(17, 2).Value = 3.674332095824184
$ws.Cells.Item(17, 3).Value = 1.055204134061853
$ws.Cells.Item(17, 4).Value = 0.1026509265787752
$ws.Cells.Item(17, 5).Value = 0.02311821252105162
$ws.Cells.Item(17, 7).Value = 0.002489596124630005
$ws.Cells.Item(17, 13).Value = 0.9143248455871742
$ws.Cells.Item(17, 14).Value = 3.336554943252764

$ws.Cells.Item(18, 2).Value = 3.609053233652617
$ws.Cells.Item(18, 3).Value = 1.033615797803463
$ws.Cells.Item(18, 4).Value = 0.1006142268706611
$ws.Cells.Item(18, 5).Value = 0.02271106404266021
$ws.Cells.Item(18, 7).Value = 0.002491686938720609
$ws.Cells.Item(18, 13).Value = 0.8974624618338822
$ws.Cells.Item(18, 14).Value = 3.310891607785436

$ws.Cells.Item(19, 2).Value = 3.586994284856246
$ws.Cells.Item(19, 3).Value = 1.026318312700482
$ws.Cells.Item(19, 4).Value = 0.09992580419061881
$ws.Cells.Item(19, 5).Value = 0.02257350122972923
$ws.Cells.Item(19, 7).Value = 0.002492399282281819
$ws.Cells.Item(19, 13).Value = 0.891763952811985
$ws.Cells.Item(19, 14).Value = 3.302211386980161

$ws.Cells.Item(20, 2).Value = 3.686434460515329
$ws.Cells.Item(20, 3).Value = 1.059205371225119
$ws.Cells.Item(20, 4).Value = 0.103028433039384
$ws.Cells.Item(20, 5).Value = 0.02319370549036037
$ws.Cells.Item(20, 7).Value = 0.002489211263837092
$ws.Cells.Item(20, 13).Value = 0.9174508584419669
$ws.Cells.Item(20, 14).Value = 3.34130893707183

$ws.Cells.Item(21, 2).Value = 4.02362956571767
$ws.Cells.Item(21, 3).Value = 1.170558980494604
$ws.Cells.Item(21, 4).Value = 0.113536524541189
$ws.Cells.Item(21, 5).Value = 0.02529825456740653
$ws.Cells.Item(21, 7).Value = 0.002478806694296386
$ws.Cells.Item(21, 13).Value = 1.004526647100874
$ws.Cells.Item(21, 14).Value = 3.473316990924786

$ws.Cells.Item(22, 2).Value = 4.246397398232943
$ws.Cells.Item(22, 3).Value = 1.244005267075295
$ws.Cells.Item(22, 4).Value = 0.12046937991893
$ws.Cells.Item(22, 5).Value = 0.0266897937315882
$ws.Cells.Item(22, 7).Value = 0.002472233620587653
$ws.Cells.Item(22, 13).Value = 1.062034150326113
$ws.Cells.Item(22, 14).Value = 3.560101064979676

$ws.Cells.Item(23, 2).Value = 4.127277436432678
$ws.Cells.Item(23, 3).Value = 1.204742246745127
$ws.Cells.Item(23, 4).Value = 0.1167630327879294
$ws.Cells.Item(23, 5).Value = 0.02594558931665958
$ws.Cells.Item(23, 7).Value = 0.002475721140749749
$ws.Cells.Item(23, 13).Value = 1.031285019572749
$ws.Cells.Item(23, 14).Value = 3.513734364879383

$ws.Cells.Item(24, 2).Value = 3.6809622795829
$ws.Cells.Item(24, 3).Value = 1.057396222815157
$ws.Cells.Item(24, 4).Value = 0.1028577437942317
$ws.Cells.Item(24, 5).Value = 0.02315957036670824
$ws.Cells.Item(24, 7).Value = 0.002489385176310189
$ws.Cells.Item(24, 13).Value = 0.9160374138105709
$ws.Cells.Item(24, 14).Value = 3.339159527965109

$ws.Cells.Item(25, 2).Value = 3.20923651614936
$ws.Cells.Item(25, 3).Value = 0.9011373680684187
$ws.Cells.Item(25, 4).Value = 0.08812033610021786
$ws.Cells.Item(25, 5).Value = 0.02021943726730058
$ws.Cells.Item(25, 7).Value = 0.002505109310176496
$ws.Cells.Item(25, 13).Value = 0.7941418404794547
$ws.Cells.Item(25, 14).Value = 3.152857768936428
